$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right after
#    the H1 title.
# ---------------------------------------------------------------------------
$metaOld = "Meta description: Read our review of Frost Queen Jackpots, with innovative graphics, bonuses, jackpots, and a chance to win up to 1,265 times your total bet. Play for free."
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $metaOld) {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Locate the closing "Prompt: ..." paragraph (now the very last paragraph)
#    and insert a brand-new bold paragraph right before it, reusing the same
#    title text that appears in the H1 heading.
# ---------------------------------------------------------------------------
$promptOld = "Prompt: Create a feature image for Frost Queen Jackpots that features a happy Maya warrior with glasses in cartoon style. The image should showcase the Ice Queen's palace in the background, with the warrior standing in front of it, holding up a pair of playing cards. The cards should have the symbols of the Ice Queen and the Maya warrior, representing the theme of the slot game. The overall tone of the image should be bright and fun, in line with the playful and exciting nature of the game."

$promptParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq $promptOld) {
        $promptParaIndex = $i
        break
    }
}

$titleText = "Play Frost Queen Jackpots for Free " + [char]0x2013 + " Innovative Graphics and Bonuses"

$promptPara = $d.Paragraphs($promptParaIndex)
$insertionPoint = $d.Range($promptPara.Range.Start, $promptPara.Range.Start)

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $titleText + '</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($xmlFrag)

# The fragment above leaves one now-empty paragraph behind (the merge target
# for the trailing "<w:p/>"); remove it so the "Prompt" paragraph goes right
# back to being the paragraph that directly follows our new bold paragraph.
$newTitleParaIndex = $promptParaIndex
$emptyParaIndex = $newTitleParaIndex + 1
if ($d.Paragraphs($emptyParaIndex).Range.Text.TrimEnd([char]13, [char]7) -eq "") {
    $d.Paragraphs($emptyParaIndex).Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Swap the old "Prompt: ..." copy for the new meta-description sentence,
#    keeping the paragraph's existing (italic) run formatting untouched.
# ---------------------------------------------------------------------------
$newDescription = "Read our review of Frost Queen Jackpots, with innovative graphics, bonuses, jackpots, and a chance to win up to 1,265 times your total bet. Play for free."
$d.Content.Find.Execute($promptOld, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)
